# Update omni-schema so that it parses Benchmark-001.yaml example
$wb = $excel.ActiveWorkbook

# --- IOFile sheet: name,path -> path,id,name,description ---
$ioFile = $wb.Worksheets.Item("IOFile")
$ioFile.Range("A1").Value = "path"
$ioFile.Range("B1").Value = "id"
$ioFile.Range("C1").Value = "name"
$ioFile.Range("D1").Value = "description"

# --- Parameter sheet: name -> values ---
$parameter = $wb.Worksheets.Item("Parameter")
$parameter.Range("A1").Value = "values"
